# Apply the "Add files via upload" commit: populate additional typing-test
# results (students #2..#21) on sheet "314" and set the print setup for
# sheet "315".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("314")
$ws2 = $wb.Worksheets.Item("315")

# --- Sheet "314": fill in the D (date) style for every new data row first,
# by copying the format from the already-styled D2 cell, so the new date
# values inherit the existing "mm\"월\" dd\"일\"" number format (style index 1)
# instead of creating a brand-new style entry.
$ws1.Range("D2").Copy() | Out-Null
$ws1.Range("D3:D21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row-by-row data -------------------------------------------------
# Column A = pass/fail ("ok"), Column C = name, Column D = date, Column E = typed count.
# Rows 2,4,7,8,16,17,18,22 already had name/date/count; they only gain the "ok" flag.
# Rows 3,5,6,9,10,11,12,13,14,15,19,20,21 are brand-new rows of data.

$ws1.Range("A2").Value2 = "ok"

$ws1.Range("A3").Value2 = "ok"
$ws1.Range("C3").Value2 = "강종희"
$ws1.Range("D3").Value2 = 43637
$ws1.Range("E3").Value2 = 210

$ws1.Range("A4").Value2 = "ok"

$ws1.Range("C5").Value2 = "김명석"
$ws1.Range("D5").Value2 = 43637
$ws1.Range("E5").Value2 = 168

$ws1.Range("C6").Value2 = "김성진 "
$ws1.Range("D6").Value2 = 43637
$ws1.Range("E6").Value2 = 172

$ws1.Range("A7").Value2 = "ok"

$ws1.Range("A8").Value2 = "ok"

$ws1.Range("A9").Value2 = "ok"
$ws1.Range("C9").Value2 = "김형민"
$ws1.Range("D9").Value2 = 43637
$ws1.Range("E9").Value2 = 218

$ws1.Range("A10").Value2 = "ok"
$ws1.Range("C10").Value2 = "박재형"
$ws1.Range("D10").Value2 = 43637
$ws1.Range("E10").Value2 = 213

$ws1.Range("A11").Value2 = "ok"
$ws1.Range("C11").Value2 = "박형준"
$ws1.Range("D11").Value2 = 43637
$ws1.Range("E11").Value2 = 237

$ws1.Range("A12").Value2 = "ok"
$ws1.Range("C12").Value2 = "백승욱 "
$ws1.Range("D12").Value2 = 43637
$ws1.Range("E12").Value2 = 227

$ws1.Range("C13").Value2 = "서찬우"
$ws1.Range("D13").Value2 = 43637
$ws1.Range("E13").Value2 = 195

$ws1.Range("C14").Value2 = "신광민"
$ws1.Range("D14").Value2 = 43637
$ws1.Range("E14").Value2 = 195

$ws1.Range("A15").Value2 = "ok"
$ws1.Range("C15").Value2 = "신용훈"
$ws1.Range("D15").Value2 = 43637
$ws1.Range("E15").Value2 = 201

$ws1.Range("A16").Value2 = "ok"

$ws1.Range("A17").Value2 = "ok"

$ws1.Range("A18").Value2 = "ok"

$ws1.Range("A19").Value2 = "ok"
$ws1.Range("C19").Value2 = "임영진"
$ws1.Range("D19").Value2 = 43637
$ws1.Range("E19").Value2 = 200

$ws1.Range("C20").Value2 = "전병현"
$ws1.Range("D20").Value2 = 43637
$ws1.Range("E20").Value2 = 134

$ws1.Range("C21").Value2 = "정행곤"
$ws1.Range("D21").Value2 = 43637
$ws1.Range("E21").Value2 = 151

$ws1.Range("A22").Value2 = "ok"

# Move the active selection to A23, matching the saved workbook state.
$ws1.Range("A23").Select() | Out-Null

# --- Sheet "315": set the page setup (paper size 9 = A4, landscape) -----
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 2

# Restore sheet "315" as the active/visible tab (selecting on sheet "314"
# above switches the active sheet, but the workbook was saved with "315"
# as the active tab).
$ws2.Activate() | Out-Null
